# Weekly update: insert two new rows of price data (rows 261-262) ahead of
# the existing series, pushing the rest of the historical rows down by two
# (261->263 ... 358->360).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 261; everything from 261 downward shifts to 263+.
$ws.Rows("261:262").Insert()

# Row 261 - new weekly record (numeric serial date 44924 = 2022-12-29)
$ws.Cells.Item(261, 1).Value = 7
$ws.Cells.Item(261, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(261, 3).Value = "Ñuble"
$ws.Cells.Item(261, 4).Value = 44924
$ws.Cells.Item(261, 5).Value = 16
$ws.Cells.Item(261, 6).Value = 100114013
$ws.Cells.Item(261, 7).Value = "Zanahoria"
$ws.Cells.Item(261, 8).Value = "Sin especificar"
$ws.Cells.Item(261, 9).Value = "Primera"
$ws.Cells.Item(261, 10).Value = 160
$ws.Cells.Item(261, 11).Value = 8000
$ws.Cells.Item(261, 12).Value = 8500
$ws.Cells.Item(261, 13).Value = 8250
$ws.Cells.Item(261, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(261, 15).Value = "Región de Ñuble"
$ws.Cells.Item(261, 16).Value = 412
$ws.Cells.Item(261, 17).Value = 20
$ws.Cells.Item(261, 18).Value = "Hortaliza"

# Row 262 - new weekly record (numeric serial date 44924 = 2022-12-29)
$ws.Cells.Item(262, 1).Value = 7
$ws.Cells.Item(262, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(262, 3).Value = "Ñuble"
$ws.Cells.Item(262, 4).Value = 44924
$ws.Cells.Item(262, 5).Value = 16
$ws.Cells.Item(262, 6).Value = 100114013
$ws.Cells.Item(262, 7).Value = "Zanahoria"
$ws.Cells.Item(262, 8).Value = "Sin especificar"
$ws.Cells.Item(262, 9).Value = "Segunda"
$ws.Cells.Item(262, 10).Value = 200
$ws.Cells.Item(262, 11).Value = 7000
$ws.Cells.Item(262, 12).Value = 7500
$ws.Cells.Item(262, 13).Value = 7250
$ws.Cells.Item(262, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(262, 15).Value = "Región de Ñuble"
$ws.Cells.Item(262, 16).Value = 362
$ws.Cells.Item(262, 17).Value = 20
$ws.Cells.Item(262, 18).Value = "Hortaliza"
